$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 160-161; existing rows 160-248 shift down to 162-250
$ws.Rows("160:161").Insert()

$ws.Range("A160").Value = 11
$ws.Range("B160").Value = "Vega Monumental Concepción"
$ws.Range("C160").Value = "Bíobío"
$ws.Range("D160").Value = 44567
$ws.Range("E160").Value = 8
$ws.Range("F160").Value = "Fruta"
$ws.Range("G160").Value = 100101
$ws.Range("H160").Value = "Berries"
$ws.Range("I160").Value = 100112025
$ws.Range("J160").Value = "Frutilla"
$ws.Range("K160").Value = "Sin especificar"
$ws.Range("L160").Value = "Primera"
$ws.Range("M160").Value = 250
$ws.Range("N160").Value = 7000
$ws.Range("O160").Value = 7500
$ws.Range("P160").Value = 7260
$ws.Range("Q160").Value = "`$/caja 7 kilos"
$ws.Range("R160").Value = "Región del Maule"
$ws.Range("S160").Value = 1037
$ws.Range("T160").Value = 7

$ws.Range("A161").Value = 11
$ws.Range("B161").Value = "Vega Monumental Concepción"
$ws.Range("C161").Value = "Bíobío"
$ws.Range("D161").Value = 44567
$ws.Range("E161").Value = 8
$ws.Range("F161").Value = "Fruta"
$ws.Range("G161").Value = 100101
$ws.Range("H161").Value = "Berries"
$ws.Range("I161").Value = 100112025
$ws.Range("J161").Value = "Frutilla"
$ws.Range("K161").Value = "Sin especificar"
$ws.Range("L161").Value = "Segunda"
$ws.Range("M161").Value = 270
$ws.Range("N161").Value = 6000
$ws.Range("O161").Value = 6500
$ws.Range("P161").Value = 6259
$ws.Range("Q161").Value = "`$/caja 7 kilos"
$ws.Range("R161").Value = "Región del Maule"
$ws.Range("S161").Value = 894
$ws.Range("T161").Value = 7

